# updated variables on participant-level page
# Adds four new "KDE" rows (QOL_KDE_Text, QOL_KDE_Sim, VOL_KDE_Text, VOL_KDE_Sim),
# each inserted immediately above its corresponding existing "KDMA" row, and
# switches the existing KDMA rows' "Labels" column from "Hyperlink to Graph" to
# "Number" (the KDE row now owns the hyperlink-to-graph label).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert blank rows from the bottom up so earlier row numbers stay stable
# while we work. Before any inserts:
#   18 QOL_KDMA_Text   19 QOL_KDMA_Sim   20 VOL_KDMA_Text   21 VOL_KDMA_Sim (style s=9)
$ws.Rows("21").Insert()
$ws.Rows("20").Insert()
$ws.Rows("19").Insert()
$ws.Rows("18").Insert()

# After the four inserts the sheet looks like:
#   18 <blank>  19 QOL_KDMA_Text  20 <blank>  21 QOL_KDMA_Sim
#   22 <blank>  23 VOL_KDMA_Text  24 <blank>  25 VOL_KDMA_Sim (style s=9)

# New row 18: QOL_KDE_Text (copy of QOL_KDMA_Text row, renamed, keeps hyperlink label)
$ws.Range("A18").Value = "Attribute Assessment"
$ws.Range("B18").Value = "QOL_KDE_Text"
$ws.Range("C18").Value = "KDMA measurement from text probe responses for the ST Quality of Life scenario"
$ws.Range("D18").Value = "Hyperlink to Graph"
$ws.Range("E18").Value = "-"
$ws.Range("F18").Value = "From TA1 Server"

# Existing row 19: QOL_KDMA_Text now reports a Number instead of the hyperlink
$ws.Range("D19").Value = "Number"

# New row 20: QOL_KDE_Sim (copy of QOL_KDMA_Sim row, renamed, keeps hyperlink label)
$ws.Range("A20").Value = "Attribute Assessment"
$ws.Range("B20").Value = "QOL_KDE_Sim"
$ws.Range("C20").Value = "KDMA measurement from sim probe responses for the ST Quality of Life Scenario"
$ws.Range("D20").Value = "Hyperlink to Graph"
$ws.Range("E20").Value = "-"
$ws.Range("F20").Value = "From TA1 Server"

# Existing row 21: QOL_KDMA_Sim now reports a Number instead of the hyperlink
$ws.Range("D21").Value = "Number"

# New row 22: VOL_KDE_Text (copy of VOL_KDMA_Text row, renamed, keeps hyperlink label)
$ws.Range("A22").Value = "Attribute Assessment"
$ws.Range("B22").Value = "VOL_KDE_Text"
$ws.Range("C22").Value = "KDMA measurement from text probe responses for the ST Value of Life scenario"
$ws.Range("D22").Value = "Hyperlink to Graph"
$ws.Range("E22").Value = "-"
$ws.Range("F22").Value = "From TA1 Server"

# Existing row 23: VOL_KDMA_Text now reports a Number instead of the hyperlink
$ws.Range("D23").Value = "Number"

# New row 24: VOL_KDE_Sim (copy of VOL_KDMA_Sim row, renamed, keeps hyperlink label).
# The original VOL_KDMA_Sim cell (now row 25) carried the special italic/custom
# font style in column C; mirror it onto the new row's C cell too.
$ws.Range("A24").Value = "Attribute Assessment"
$ws.Range("B24").Value = "VOL_KDE_Sim"
$ws.Range("C24").Value = "KDMA measurement from sim probe responses for the ST Value of Life Scenario"
$ws.Range("D24").Value = "Hyperlink to Graph"
$ws.Range("E24").Value = "-"
$ws.Range("F24").Value = "From TA1 Server"
$ws.Range("C25").Copy()
$ws.Range("C24").PasteSpecial(-4122)

# Existing row 25: VOL_KDMA_Sim now reports a Number instead of the hyperlink
$ws.Range("D25").Value = "Number"

# Sheet view / selection housekeeping to match the saved workbook state
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("D25").Select()
